$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (Changed) date column C for rows 2 through 45
# from serial date 45741 (2025-03-25) to 45742 (2025-03-26)
foreach ($row in 2..45) {
    $ws.Cells.Item($row, 3).Value = 45742
}
